# Apply scraped Leve profit-sheet value updates (scheduled runner refresh).
# Each write targets plain data cells (no formulas in this workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 294.74075  # H33: 277.5 -> 294.74075
$ws.Cells.Item(33, 9).Value = 228.09525  # I33: 208.57143 -> 228.09525
$ws.Cells.Item(33, 10).Value = 528  # J33: 567 -> 528
$ws.Cells.Item(33, 11).Value = 228.09525  # K33: 208.57143 -> 228.09525
$ws.Cells.Item(33, 12).Value = 528  # L33: 567 -> 528
$ws.Cells.Item(33, 13).Value = 0.904750000000007  # M33: 20.42857000000001 -> 0.904750000000007
$ws.Cells.Item(33, 14).Value = -986  # N33: -1025 -> -986

$ws.Cells.Item(113, 8).Value = 2100  # H113: 2006 -> 2100
$ws.Cells.Item(113, 10).Value = 2100  # J113: 2006 -> 2100
$ws.Cells.Item(113, 12).Value = 2100  # L113: 2006 -> 2100
$ws.Cells.Item(113, 14).Value = -8608  # N113: -8514 -> -8608

$ws.Cells.Item(135, 8).Value = 33333928  # H135: 32258644 -> 33333928
$ws.Cells.Item(135, 9).Value = 364.85184  # I135: 355.39285 -> 364.85184
$ws.Cells.Item(135, 11).Value = 3283.66656  # K135: 3198.53565 -> 3283.66656
$ws.Cells.Item(135, 13).Value = -748.6665599999997  # M135: -663.5356500000003 -> -748.6665599999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 877.4  # H2: 926.5714 -> 877.4
$ws.Cells.Item(2, 9).Value = 686.375  # I2: 757.4286 -> 686.375
$ws.Cells.Item(2, 11).Value = 686.375  # K2: 757.4286 -> 686.375
$ws.Cells.Item(2, 13).Value = -573.375  # M2: -644.4286 -> -573.375

$ws.Cells.Item(49, 8).Value = 10000  # H49: 9500 -> 10000
$ws.Cells.Item(49, 10).Value = 10000  # J49: 9500 -> 10000
$ws.Cells.Item(49, 12).Value = 10000  # L49: 9500 -> 10000
$ws.Cells.Item(49, 14).Value = -10520  # N49: -10020 -> -10520

$ws.Cells.Item(61, 8).Value = 76924216  # H61: 83334610 -> 76924216
$ws.Cells.Item(61, 9).Value = 100000800  # I61: 111112056 -> 100000800
$ws.Cells.Item(61, 11).Value = 100000800  # K61: 111112056 -> 100000800
$ws.Cells.Item(61, 13).Value = -100000588  # M61: -111111844 -> -100000588

$ws.Cells.Item(116, 8).Value = 877.4  # H116: 926.5714 -> 877.4
$ws.Cells.Item(116, 9).Value = 686.375  # I116: 757.4286 -> 686.375
$ws.Cells.Item(116, 11).Value = 686.375  # K116: 757.4286 -> 686.375
$ws.Cells.Item(116, 13).Value = 1607.625  # M116: 1536.5714 -> 1607.625

$ws.Cells.Item(122, 8).Value = 2256  # H122: 2885.2 -> 2256
$ws.Cells.Item(122, 9).Value = 1970  # I122: 2637.3333 -> 1970
$ws.Cells.Item(122, 11).Value = 5910  # K122: 7911.999899999999 -> 5910
$ws.Cells.Item(122, 13).Value = -3460  # M122: -5461.999899999999 -> -3460

$ws.Cells.Item(136, 8).Value = 76924216  # H136: 83334610 -> 76924216
$ws.Cells.Item(136, 9).Value = 100000800  # I136: 111112056 -> 100000800
$ws.Cells.Item(136, 11).Value = 300002400  # K136: 333336168 -> 300002400
$ws.Cells.Item(136, 13).Value = -299999850  # M136: -333333618 -> -299999850

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 877.4  # H3: 926.5714 -> 877.4
$ws.Cells.Item(3, 9).Value = 686.375  # I3: 757.4286 -> 686.375
$ws.Cells.Item(3, 11).Value = 686.375  # K3: 757.4286 -> 686.375
$ws.Cells.Item(3, 13).Value = -572.375  # M3: -643.4286 -> -572.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(44, 8).Value = 0  # H44: 3500 -> 0
$ws.Cells.Item(44, 10).Value = 0  # J44: 3500 -> 0
$ws.Cells.Item(44, 12).Value = 0  # L44: 3500 -> 0
$ws.Cells.Item(44, 14).Value = $null  # N44: -4384 -> (empty)

$ws.Cells.Item(86, 8).Value = 4793788  # H86: 5592011 -> 4793788
$ws.Cells.Item(86, 10).Value = 36935.555  # J86: 46217.715 -> 36935.555
$ws.Cells.Item(86, 12).Value = 36935.555  # L86: 46217.715 -> 36935.555
$ws.Cells.Item(86, 14).Value = -39181.555  # N86: -48463.715 -> -39181.555

$ws.Cells.Item(89, 8).Value = 4793788  # H89: 5592011 -> 4793788
$ws.Cells.Item(89, 10).Value = 36935.555  # J89: 46217.715 -> 36935.555
$ws.Cells.Item(89, 12).Value = 184677.775  # L89: 231088.575 -> 184677.775
$ws.Cells.Item(89, 14).Value = -195909.775  # N89: -242320.575 -> -195909.775

$ws.Cells.Item(99, 8).Value = 1523.8182  # H99: 1528.3636 -> 1523.8182
$ws.Cells.Item(99, 9).Value = 1541.6666  # I99: 1464.2858 -> 1541.6666
$ws.Cells.Item(99, 10).Value = 1502.4  # J99: 1640.5 -> 1502.4
$ws.Cells.Item(99, 11).Value = 1541.6666  # K99: 1464.2858 -> 1541.6666
$ws.Cells.Item(99, 12).Value = 1502.4  # L99: 1640.5 -> 1502.4
$ws.Cells.Item(99, 13).Value = -43.66660000000002  # M99: 33.71419999999989 -> -43.66660000000002
$ws.Cells.Item(99, 14).Value = -4498.4  # N99: -4636.5 -> -4498.4

$ws.Cells.Item(122, 8).Value = 757.3333  # H122: 650.2 -> 757.3333
$ws.Cells.Item(122, 9).Value = 775.61536  # I122: 666.2941 -> 775.61536
$ws.Cells.Item(122, 10).Value = 638.5  # J122: 559 -> 638.5
$ws.Cells.Item(122, 11).Value = 2326.84608  # K122: 1998.8823 -> 2326.84608
$ws.Cells.Item(122, 12).Value = 1915.5  # L122: 1677 -> 1915.5
$ws.Cells.Item(122, 13).Value = 123.1539199999997  # M122: 451.1177000000002 -> 123.1539199999997
$ws.Cells.Item(122, 14).Value = -6815.5  # N122: -6577 -> -6815.5

$ws.Cells.Item(126, 8).Value = 1523.8182  # H126: 1528.3636 -> 1523.8182
$ws.Cells.Item(126, 9).Value = 1541.6666  # I126: 1464.2858 -> 1541.6666
$ws.Cells.Item(126, 10).Value = 1502.4  # J126: 1640.5 -> 1502.4
$ws.Cells.Item(126, 11).Value = 4624.9998  # K126: 4392.857400000001 -> 4624.9998
$ws.Cells.Item(126, 12).Value = 4507.200000000001  # L126: 4921.5 -> 4507.200000000001
$ws.Cells.Item(126, 13).Value = -2154.9998  # M126: -1922.857400000001 -> -2154.9998
$ws.Cells.Item(126, 14).Value = -9447.200000000001  # N126: -9861.5 -> -9447.200000000001

$ws.Cells.Item(132, 8).Value = 1823.2174  # H132: 1824.5217 -> 1823.2174
$ws.Cells.Item(132, 9).Value = 1296.2858  # I132: 1298.4286 -> 1296.2858
$ws.Cells.Item(132, 11).Value = 3888.8574  # K132: 3895.2858 -> 3888.8574
$ws.Cells.Item(132, 13).Value = -1358.8574  # M132: -1365.2858 -> -1358.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 300  # H25: 2000 -> 300
$ws.Cells.Item(25, 10).Value = 300  # J25: 2000 -> 300
$ws.Cells.Item(25, 12).Value = 900  # L25: 6000 -> 900
$ws.Cells.Item(25, 14).Value = -1238  # N25: -6338 -> -1238

$ws.Cells.Item(26, 8).Value = 250  # H26: 184.85715 -> 250
$ws.Cells.Item(26, 9).Value = 220  # I26: 148.8 -> 220
$ws.Cells.Item(26, 10).Value = 260  # J26: 275 -> 260
$ws.Cells.Item(26, 11).Value = 660  # K26: 446.4 -> 660
$ws.Cells.Item(26, 12).Value = 780  # L26: 825 -> 780
$ws.Cells.Item(26, 13).Value = -372  # M26: -158.4 -> -372
$ws.Cells.Item(26, 14).Value = -1356  # N26: -1401 -> -1356

$ws.Cells.Item(29, 8).Value = 665.6667  # H29: 700.4 -> 665.6667
$ws.Cells.Item(29, 9).Value = 296  # I29: 100 -> 296
$ws.Cells.Item(29, 11).Value = 888  # K29: 300 -> 888
$ws.Cells.Item(29, 13).Value = -611  # M29: -23 -> -611

$ws.Cells.Item(30, 8).Value = 300  # H30: 2000 -> 300
$ws.Cells.Item(30, 10).Value = 300  # J30: 2000 -> 300
$ws.Cells.Item(30, 12).Value = 900  # L30: 6000 -> 900
$ws.Cells.Item(30, 14).Value = -1104  # N30: -6204 -> -1104

$ws.Cells.Item(31, 8).Value = 1520.6  # H31: 1500.625 -> 1520.6
$ws.Cells.Item(31, 9).Value = 533.6667  # I31: 600.2 -> 533.6667
$ws.Cells.Item(31, 10).Value = 3001  # J31: 3001.3333 -> 3001
$ws.Cells.Item(31, 11).Value = 1601.0001  # K31: 1800.6 -> 1601.0001
$ws.Cells.Item(31, 12).Value = 9003  # L31: 9003.999899999999 -> 9003
$ws.Cells.Item(31, 13).Value = -1313.0001  # M31: -1512.6 -> -1313.0001
$ws.Cells.Item(31, 14).Value = -9579  # N31: -9579.999899999999 -> -9579

$ws.Cells.Item(32, 8).Value = 1911.3334  # H32: 551 -> 1911.3334
$ws.Cells.Item(32, 10).Value = 2300  # J32: 0 -> 2300
$ws.Cells.Item(32, 12).Value = 6900  # L32: 0 -> 6900
$ws.Cells.Item(32, 14).Value = -7466  # N32: (empty) -> -7466

$ws.Cells.Item(33, 8).Value = 301.18182  # H33: 245.85715 -> 301.18182
$ws.Cells.Item(33, 9).Value = 135.33333  # I33: 170 -> 135.33333
$ws.Cells.Item(33, 10).Value = 500.2  # J33: 524 -> 500.2
$ws.Cells.Item(33, 11).Value = 811.9999799999999  # K33: 1020 -> 811.9999799999999
$ws.Cells.Item(33, 12).Value = 3001.2  # L33: 3144 -> 3001.2
$ws.Cells.Item(33, 13).Value = -528.9999799999999  # M33: -737 -> -528.9999799999999
$ws.Cells.Item(33, 14).Value = -3567.2  # N33: -3710 -> -3567.2

$ws.Cells.Item(34, 8).Value = 8334939.5  # H34: 10001907 -> 8334939.5
$ws.Cells.Item(34, 9).Value = 173  # I34: 246 -> 173
$ws.Cells.Item(34, 11).Value = 519  # K34: 738 -> 519
$ws.Cells.Item(34, 13).Value = -435  # M34: -654 -> -435

$ws.Cells.Item(35, 8).Value = 3401.6667  # H35: 4001.6667 -> 3401.6667
$ws.Cells.Item(35, 9).Value = 200  # I35: 0 -> 200
$ws.Cells.Item(35, 10).Value = 5002.5  # J35: 4001.6667 -> 5002.5
$ws.Cells.Item(35, 11).Value = 600  # K35: 0 -> 600
$ws.Cells.Item(35, 12).Value = 15007.5  # L35: 12005.0001 -> 15007.5
$ws.Cells.Item(35, 13).Value = -312  # M35: (empty) -> -312
$ws.Cells.Item(35, 14).Value = -15583.5  # N35: -12581.0001 -> -15583.5

$ws.Cells.Item(36, 8).Value = 246  # H36: 0 -> 246
$ws.Cells.Item(36, 9).Value = 402  # I36: 0 -> 402
$ws.Cells.Item(36, 10).Value = 90  # J36: 0 -> 90
$ws.Cells.Item(36, 11).Value = 1206  # K36: 0 -> 1206
$ws.Cells.Item(36, 12).Value = 270  # L36: 0 -> 270
$ws.Cells.Item(36, 13).Value = -1037  # M36: (empty) -> -1037
$ws.Cells.Item(36, 14).Value = -608  # N36: (empty) -> -608

$ws.Cells.Item(38, 8).Value = 0  # H38: 110 -> 0
$ws.Cells.Item(38, 9).Value = 0  # I38: 40 -> 0
$ws.Cells.Item(38, 10).Value = 0  # J38: 250 -> 0
$ws.Cells.Item(38, 11).Value = 0  # K38: 120 -> 0
$ws.Cells.Item(38, 12).Value = 0  # L38: 750 -> 0
$ws.Cells.Item(38, 13).Value = $null  # M38: 227 -> (empty)
$ws.Cells.Item(38, 14).Value = $null  # N38: -1444 -> (empty)

$ws.Cells.Item(39, 8).Value = 2705.7307  # H39: 2915.739 -> 2705.7307
$ws.Cells.Item(39, 10).Value = 2614.5417  # J39: 2831.524 -> 2614.5417
$ws.Cells.Item(39, 12).Value = 7843.625100000001  # L39: 8494.572 -> 7843.625100000001
$ws.Cells.Item(39, 14).Value = -8431.625100000001  # N39: -9082.572 -> -8431.625100000001

$ws.Cells.Item(41, 8).Value = 507.69232  # H41: 638.46155 -> 507.69232
$ws.Cells.Item(41, 9).Value = 300  # I41: 0 -> 300
$ws.Cells.Item(41, 10).Value = 525  # J41: 638.46155 -> 525
$ws.Cells.Item(41, 11).Value = 900  # K41: 0 -> 900
$ws.Cells.Item(41, 12).Value = 1575  # L41: 1915.38465 -> 1575
$ws.Cells.Item(41, 13).Value = -562  # M41: (empty) -> -562
$ws.Cells.Item(41, 14).Value = -2251  # N41: -2591.38465 -> -2251

$ws.Cells.Item(42, 8).Value = 4253.3335  # H42: 4503.2 -> 4253.3335
$ws.Cells.Item(42, 10).Value = 4253.3335  # J42: 4503.2 -> 4253.3335
$ws.Cells.Item(42, 12).Value = 12760.0005  # L42: 13509.6 -> 12760.0005
$ws.Cells.Item(42, 14).Value = -13828.0005  # N42: -14577.6 -> -13828.0005

$ws.Cells.Item(43, 8).Value = 900  # H43: 0 -> 900
$ws.Cells.Item(43, 9).Value = 900  # I43: 0 -> 900
$ws.Cells.Item(43, 11).Value = 2700  # K43: 0 -> 2700
$ws.Cells.Item(43, 13).Value = -2586  # M43: (empty) -> -2586

$ws.Cells.Item(44, 8).Value = 1019.25  # H44: 679.3 -> 1019.25
$ws.Cells.Item(44, 9).Value = 262.5  # I44: 298.6 -> 262.5
$ws.Cells.Item(44, 10).Value = 1776  # J44: 1060 -> 1776
$ws.Cells.Item(44, 11).Value = 787.5  # K44: 895.8000000000001 -> 787.5
$ws.Cells.Item(44, 12).Value = 5328  # L44: 3180 -> 5328
$ws.Cells.Item(44, 13).Value = -389.5  # M44: -497.8000000000001 -> -389.5
$ws.Cells.Item(44, 14).Value = -6124  # N44: -3976 -> -6124

$ws.Cells.Item(46, 8).Value = 461.8  # H46: 1000 -> 461.8
$ws.Cells.Item(46, 9).Value = 351.25  # I46: 0 -> 351.25
$ws.Cells.Item(46, 10).Value = 904  # J46: 1000 -> 904
$ws.Cells.Item(46, 11).Value = 1053.75  # K46: 0 -> 1053.75
$ws.Cells.Item(46, 12).Value = 2712  # L46: 3000 -> 2712
$ws.Cells.Item(46, 13).Value = -962.75  # M46: (empty) -> -962.75
$ws.Cells.Item(46, 14).Value = -2894  # N46: -3182 -> -2894

$ws.Cells.Item(47, 8).Value = 60  # H47: 2335.7144 -> 60
$ws.Cells.Item(47, 9).Value = 60  # I47: 87.5 -> 60
$ws.Cells.Item(47, 10).Value = 0  # J47: 5333.3335 -> 0
$ws.Cells.Item(47, 11).Value = 180  # K47: 262.5 -> 180
$ws.Cells.Item(47, 12).Value = 0  # L47: 16000.0005 -> 0
$ws.Cells.Item(47, 13).Value = 251  # M47: 168.5 -> 251
$ws.Cells.Item(47, 14).Value = $null  # N47: -16862.0005 -> (empty)

$ws.Cells.Item(140, 8).Value = 23396.123  # H140: 23263.469 -> 23396.123
$ws.Cells.Item(140, 9).Value = 52331.35  # I140: 52006.35 -> 52331.35
$ws.Cells.Item(140, 11).Value = 156994.05  # K140: 156019.05 -> 156994.05
$ws.Cells.Item(140, 13).Value = -151814.05  # M140: -150839.05 -> -151814.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1363.3334  # H113: 1282.8572 -> 1363.3334
$ws.Cells.Item(113, 9).Value = 1237.1428  # I113: 1140 -> 1237.1428
$ws.Cells.Item(113, 11).Value = 1237.1428  # K113: 1140 -> 1237.1428
$ws.Cells.Item(113, 13).Value = 932.8571999999999  # M113: 1030 -> 932.8571999999999

$ws.Cells.Item(126, 8).Value = 2193.5715  # H126: 2161.4 -> 2193.5715
$ws.Cells.Item(126, 9).Value = 1810.25  # I126: 1800 -> 1810.25
$ws.Cells.Item(126, 10).Value = 2704.6667  # J126: 3004.6667 -> 2704.6667
$ws.Cells.Item(126, 11).Value = 5430.75  # K126: 5400 -> 5430.75
$ws.Cells.Item(126, 12).Value = 8114.000100000001  # L126: 9014.000100000001 -> 8114.000100000001
$ws.Cells.Item(126, 13).Value = -2960.75  # M126: -2930 -> -2960.75
$ws.Cells.Item(126, 14).Value = -13054.0001  # N126: -13954.0001 -> -13054.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3182.4092  # H40: 2255.92 -> 3182.4092
$ws.Cells.Item(40, 9).Value = 2212.2354  # I40: 2094.9048 -> 2212.2354
$ws.Cells.Item(40, 10).Value = 6481  # J40: 3101.25 -> 6481
$ws.Cells.Item(40, 11).Value = 2212.2354  # K40: 2094.9048 -> 2212.2354
$ws.Cells.Item(40, 12).Value = 6481  # L40: 3101.25 -> 6481
$ws.Cells.Item(40, 13).Value = -2076.2354  # M40: -1958.9048 -> -2076.2354
$ws.Cells.Item(40, 14).Value = -6753  # N40: -3373.25 -> -6753

$ws.Cells.Item(136, 8).Value = 1820.7  # H136: 1511.8889 -> 1820.7
$ws.Cells.Item(136, 9).Value = 1372.4286  # I136: 1325.875 -> 1372.4286
$ws.Cells.Item(136, 10).Value = 2866.6667  # J136: 3000 -> 2866.6667
$ws.Cells.Item(136, 11).Value = 4117.2858  # K136: 3977.625 -> 4117.2858
$ws.Cells.Item(136, 12).Value = 8600.000100000001  # L136: 9000 -> 8600.000100000001
$ws.Cells.Item(136, 13).Value = -1567.2858  # M136: -1427.625 -> -1567.2858
$ws.Cells.Item(136, 14).Value = -13700.0001  # N136: -14100 -> -13700.0001

$ws.Cells.Item(140, 8).Value = 48883  # H140: 50688 -> 48883
$ws.Cells.Item(140, 10).Value = 48883  # J140: 50688 -> 48883
$ws.Cells.Item(140, 12).Value = 48883  # L140: 50688 -> 48883
$ws.Cells.Item(140, 14).Value = -59243  # N140: -61048 -> -59243

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(98, 8).Value = 21097.8  # H98: 21750 -> 21097.8
$ws.Cells.Item(98, 10).Value = 21097.8  # J98: 21750 -> 21097.8
$ws.Cells.Item(98, 12).Value = 21097.8  # L98: 21750 -> 21097.8
$ws.Cells.Item(98, 14).Value = -27087.8  # N98: -27740 -> -27087.8

$ws.Cells.Item(107, 8).Value = 637.2727  # H107: 484.2353 -> 637.2727
$ws.Cells.Item(107, 9).Value = 488.75  # I107: 405.16666 -> 488.75
$ws.Cells.Item(107, 10).Value = 1033.3334  # J107: 674 -> 1033.3334
$ws.Cells.Item(107, 11).Value = 1466.25  # K107: 1215.49998 -> 1466.25
$ws.Cells.Item(107, 12).Value = 3100.0002  # L107: 2022 -> 3100.0002
$ws.Cells.Item(107, 13).Value = 453.75  # M107: 704.5000199999999 -> 453.75
$ws.Cells.Item(107, 14).Value = -6940.0002  # N107: -5862 -> -6940.0002

$ws.Cells.Item(136, 8).Value = 956.1111  # H136: 979.1429000000001 -> 956.1111
$ws.Cells.Item(136, 9).Value = 889.5  # I136: 919.08 -> 889.5
$ws.Cells.Item(136, 11).Value = 2668.5  # K136: 2757.24 -> 2668.5
$ws.Cells.Item(136, 13).Value = -118.5  # M136: -207.2400000000002 -> -118.5

$ws.Cells.Item(138, 8).Value = 34766.668  # H138: 34000 -> 34766.668
$ws.Cells.Item(138, 10).Value = 34766.668  # J138: 34000 -> 34766.668
$ws.Cells.Item(138, 12).Value = 34766.668  # L138: 34000 -> 34766.668
$ws.Cells.Item(138, 14).Value = -45046.668  # N138: -44280 -> -45046.668
